$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111965439
$ws.Range("B2").Value = 56430
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 734927
$ws.Range("R2").Value = 7088234
$ws.Range("Z2").Value = "16:40"
$ws.Range("AB2").Value = "16:40"
$ws.Range("A3").Value = 111964863
$ws.Range("B3").Value = 89893
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 2062
$ws.Range("F3").Value = "Ulltickeporing"
$ws.Range("G3").Value = "Skeletocutis brevispora"
$ws.Range("H3").Value = "Niemelä"
$ws.Range("Q3").Value = 734972
$ws.Range("R3").Value = 7088253
$ws.Range("Z3").Value = "16:12"
$ws.Range("AB3").Value = "16:12"
$ws.Range("A4").Value = 111964457
$ws.Range("B4").Value = 56430
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 734949
$ws.Range("R4").Value = 7088269
$ws.Range("Z4").Value = "16:01"
$ws.Range("AB4").Value = "16:01"
$ws.Range("A5").Value = 111964622
$ws.Range("B5").Value = 89993
$ws.Range("E5").Value = 1209
$ws.Range("F5").Value = "Rynkskinn"
$ws.Range("G5").Value = "Phlebia centrifuga"
$ws.Range("H5").Value = "P.Karst."
$ws.Range("A7").Value = 111964050
$ws.Range("B7").Value = 90213
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 898
$ws.Range("F7").Value = "Blackticka"
$ws.Range("G7").Value = "Steccherinum collabens"
$ws.Range("H7").Value = "(Fr.) Vesterholt"
$ws.Range("Q7").Value = 734893
$ws.Range("R7").Value = 7088355
$ws.Range("Z7").Value = "15:42"
$ws.Range("AB7").Value = "15:42"
$ws.Range("A8").Value = 111965370
$ws.Range("B8").Value = 81385
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1312
$ws.Range("F8").Value = "Gammelgransskål"
$ws.Range("G8").Value = "Pseudographis pinicola"
$ws.Range("H8").Value = "(Nyl.) Rehm"
$ws.Range("Q8").Value = 734940
$ws.Range("R8").Value = 7088232
$ws.Range("Z8").Value = "16:38"
$ws.Range("AB8").Value = "16:38"
$ws.Range("A9").Value = 111964632
$ws.Range("B9").Value = 77650
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 734972
$ws.Range("R9").Value = 7088253
$ws.Range("Z9").Value = "16:12"
$ws.Range("AB9").Value = "16:12"
$ws.Range("A10").Value = 111964175
$ws.Range("B10").Value = 89571
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = "Granticka"
$ws.Range("G10").Value = "Porodaedalea chrysoloma"
$ws.Range("H10").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q10").Value = 734896
$ws.Range("R10").Value = 7088342
$ws.Range("Z10").Value = "15:42"
$ws.Range("AB10").Value = "15:42"
$ws.Range("A11").Value = 111965883
$ws.Range("B11").Value = 55643
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 102612
$ws.Range("F11").Value = "Järpe"
$ws.Range("G11").Value = "Tetrastes bonasia"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 734847
$ws.Range("R11").Value = 7088238
$ws.Range("Z11").Value = "17:05"
$ws.Range("AB11").Value = "17:05"
$ws.Range("M9").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("M11").Value = "lockläte, övriga läten"
